# #5: cash & deposit done
#
# Fixes the "存款" (deposits) sheet: row 1 previously duplicated the first
# data row instead of holding real column headers, and the per-row records
# were missing the property_category / category / date / legislator_name /
# legislator_id / source_file / index columns that every other sheet in
# this workbook has. This adds the missing header row + the seven trailing
# metadata columns (G:M) to every data row, and slides the "total" amount
# left into F (dropping the old, now-unused USD "quantity" column).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("存款")

# ---- header row (row 1) -----------------------------------------------
# B1:F1 already exist and keep their current (bold/bordered) style; they
# just need their text corrected from stray data to real header labels.
$ws.Cells.Item(1,2).Value = "bank"
$ws.Cells.Item(1,3).Value = "deposit_type"
$ws.Cells.Item(1,4).Value = "currency"
$ws.Cells.Item(1,5).Value = "owner"
$ws.Cells.Item(1,6).Value = "total"

# G1:M1 are brand new cells - copy the header style from an existing
# header cell (E1) before writing the text so they pick up style index 1
# instead of defaulting to unstyled.
$headerStyleSrc = $ws.Cells.Item(1,5)
$newHeaders = @("property_category","category","date","legislator_name","legislator_id","source_file","index")
for ($i = 0; $i -lt $newHeaders.Length; $i++) {
    $col = 7 + $i   # G=7 .. M=13
    $cell = $ws.Cells.Item(1,$col)
    $headerStyleSrc.Copy($cell)
    $cell.Value = $newHeaders[$i]
}

# ---- data rows (rows 2-14) --------------------------------------------
# columns: A=index, B=bank, C=deposit_type, D=currency, E=owner, F=total
#          G=property_category, H=category, I=date, J=legislator_name,
#          K=legislator_id, L=source_file, M=index (duplicate of A)
$rows = @(
    @{ idx=51; bank="台北富邦商業銀行";     type="活期儲蓄存款"; cur="新臺幣"; owner="丁守中"; total=139682 },
    @{ idx=52; bank="中華郵政股份有限公司"; type="活期儲蓄存款"; cur="新臺幣"; owner="丁守中"; total=244911 },
    @{ idx=53; bank="永豐商業銀行";         type="活期儲蓄存款"; cur="新臺幣"; owner="丁守中"; total=664256 },
    @{ idx=54; bank="台北富邦商業銀行";     type="活期儲蓄存款"; cur="新臺幣"; owner="溫子苓"; total=1915713 },
    @{ idx=55; bank="中華郵政股份有限公司"; type="活期存款";     cur="新臺幣"; owner="溫子苓"; total=149924 },
    @{ idx=56; bank="台北富邦商業銀行";     type="定期存款";     cur="美金";   owner="溫子苓"; total=305126 },
    @{ idx=57; bank="台北富邦商業銀行";     type="支票存款";     cur="新臺幣"; owner="溫子苓"; total=10000 },
    @{ idx=58; bank="上海商業儲蓄銀行";     type="活期儲蓄存款"; cur="新臺幣"; owner="溫子苓"; total=733 },
    @{ idx=59; bank="國泰世華商業銀行";     type="活期存款";     cur="新臺幣"; owner="溫子苓"; total=55130 },
    @{ idx=60; bank="聯邦商業銀行";         type="活期儲蓄存款"; cur="新臺幣"; owner="溫子苓"; total=5000 },
    @{ idx=61; bank="兆豐國際商業銀行";     type="活期儲蓄存款"; cur="新臺幣"; owner="溫子苓"; total=794873 },
    @{ idx=62; bank="華南商業銀行";         type="活期儲蓄存款"; cur="新臺幣"; owner="丁守中"; total=75867 },
    @{ idx=63; bank="華南商業銀行";         type="活期儲蓄存款"; cur="新臺幣"; owner="溫子苓"; total=194103 }
)

$legislatorName = "丁守中"
$legislatorId = 515
$sourceFile = "tmpf49e1"
$propertyCategory = "deposit"
$category = "normal"
$date = "2012-04-02"

for ($r = 0; $r -lt $rows.Length; $r++) {
    $row = $rows[$r]
    $excelRow = $r + 2   # data starts at row 2

    $ws.Cells.Item($excelRow,2).Value = $row.bank
    $ws.Cells.Item($excelRow,3).Value = $row.type
    $ws.Cells.Item($excelRow,4).Value = $row.cur
    $ws.Cells.Item($excelRow,5).Value = $row.owner
    $ws.Cells.Item($excelRow,6).Value = $row.total

    $ws.Cells.Item($excelRow,7).Value = $propertyCategory
    $ws.Cells.Item($excelRow,8).Value = $category
    $ws.Cells.Item($excelRow,9).Value = $date
    $ws.Cells.Item($excelRow,10).Value = $legislatorName
    $ws.Cells.Item($excelRow,11).Value = $legislatorId
    $ws.Cells.Item($excelRow,12).Value = $sourceFile
    $ws.Cells.Item($excelRow,13).Value = $row.idx
}

Write-Host "deposit sheet rebuilt"
